$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and volume-change (E) values per row, 2023-03-16 GitHub Actions refresh
$updates = @(
    @{ Row = 2; Price = '24.929.10'; Volume = '  +1.44%  ' }
    @{ Row = 3; Price = '1.664.49'; Volume = '  -0.39%  ' }
    @{ Row = 4; Price = '1.001'; Volume = '  -0.64%  ' }
    @{ Row = 5; Price = '317.17'; Volume = '  +3.20%  ' }
    @{ Row = 6; Price = '0.9984'; Volume = '  -0.52%  ' }
    @{ Row = 7; Price = '0.3639'; Volume = '  -1.19%  ' }
    @{ Row = 8; Price = '47.34'; Volume = '  -1.33%  ' }
    @{ Row = 9; Price = '0.3288'; Volume = '  -2.23%  ' }
    @{ Row = 10; Price = '1.145'; Volume = '  -2.24%  ' }
    @{ Row = 11; Price = '0.07107'; Volume = '  -2.67%  ' }
    @{ Row = 12; Price = '0.9984'; Volume = '  -0.68%  ' }
    @{ Row = 13; Price = '6.086'; Volume = '  -1.30%  ' }
    @{ Row = 14; Price = '19.75'; Volume = '  -3.36%  ' }
    @{ Row = 15; Price = '1.669.18'; Volume = '  -0.14%  ' }
    @{ Row = 16; Price = '6.654'; Volume = '  -1.94%  ' }
    @{ Row = 17; Price = '0.00001055'; Volume = '  -3.54%  ' }
    @{ Row = 18; Price = '0.06659'; Volume = '  +0.30%  ' }
    @{ Row = 19; Price = '0.9979'; Volume = '  -0.46%  ' }
    @{ Row = 20; Price = '79.83'; Volume = '  -1.90%  ' }
    @{ Row = 21; Price = '5.952'; Volume = '  -3.85%  ' }
    @{ Row = 22; Price = '15.86'; Volume = '  -5.24%  ' }
    @{ Row = 23; Price = '12.69'; Volume = '  +0.30%  ' }
    @{ Row = 24; Price = '24.893.51'; Volume = '  +1.38%  ' }
    @{ Row = 25; Price = '2.429'; Volume = '  -0.45%  ' }
    @{ Row = 26; Price = '2.427'; Volume = '  -9.51%  ' }
    @{ Row = 27; Price = '148.95'; Volume = '  +0.28%  ' }
    @{ Row = 28; Price = '18.70'; Volume = '  -5.62%  ' }
    @{ Row = 29; Price = '1.238'; Volume = '  +2.43%  ' }
    @{ Row = 30; Price = '1.851.79'; Volume = '  -0.40%  ' }
    @{ Row = 31; Price = '126.04'; Volume = '  -2.98%  ' }
    @{ Row = 32; Price = '4.124'; Volume = '  -0.16%  ' }
    @{ Row = 33; Price = '5.898'; Volume = '  -9.11%  ' }
    @{ Row = 34; Price = '0.08528'; Volume = '  -0.79%  ' }
    @{ Row = 35; Price = '1.671'; Volume = '  -2.17%  ' }
    @{ Row = 36; Price = '12.38'; Volume = '  -6.53%  ' }
    @{ Row = 37; Price = '1.284'; Volume = '  +4.54%  ' }
    @{ Row = 38; Price = '5.249'; Volume = '  -2.75%  ' }
    @{ Row = 39; Price = '0.02282'; Volume = '  -2.27%  ' }
    @{ Row = 40; Price = '0.06113'; Volume = '  -5.17%  ' }
    @{ Row = 41; Price = '8.378'; Volume = '  -4.89%  ' }
    @{ Row = 42; Price = '0.2084'; Volume = '  -3.33%  ' }
    @{ Row = 43; Price = '0.9980'; Volume = '  -0.57%  ' }
    @{ Row = 44; Price = '0.5974'; Volume = '  -4.02%  ' }
    @{ Row = 45; Price = '3.823'; Volume = '  +1.11%  ' }
    @{ Row = 46; Price = '12.88'; Volume = '  -4.10%  ' }
    @{ Row = 47; Price = '0.5677'; Volume = '  -3.88%  ' }
    @{ Row = 48; Price = '126.33'; Volume = '  +0.51%  ' }
    @{ Row = 49; Price = '1.969'; Volume = '  -3.45%  ' }
    @{ Row = 50; Price = '0.07040'; Volume = '  -1.08%  ' }
    @{ Row = 51; Price = '1.201'; Volume = '  +0.95%  ' }
)

foreach ($u in $updates) {
    $priceCell = $ws.Range("D" + $u.Row)
    # Prefix with a leading single-quote so Excel keeps numeric-looking text
    # (e.g. 24.929.10) stored as a string instead of silently coercing it to a number.
    $priceCell.Value = "'" + $u.Price
    $priceCell.Style = "Normal"

    $ws.Range("E" + $u.Row).Value = $u.Volume
}
